$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(6)

# Fix typo "Sofware" -> "Software" in the big red-bordered callout box.
# Retype the misspelled word (+ the following space) so the text reads
# "Software bugs ..." instead of "Sofware bugs ...". The box auto-fits
# ("spAutoFit"), so its height updates to match the corrected text.
$tr = $sh.TextFrame.TextRange
$span = $tr.Characters(1, 8)
[void]$span.Delete()
[void]$tr.Characters(1, 0).InsertBefore("Software ")

# The box also widens slightly to fit the corrected line.
$sh.Width = 651.9401
